$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.694.95"
$ws.Range("E2").Value = "  +1.99%  "

$ws.Range("D3").Value = "1.565.37"
$ws.Range("E3").Value = "  -0.19%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.990"
$ws.Range("E4").Value = "  -1.89%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.67"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.493"
$ws.Range("E6").Value = "  +0.40%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.989"
$ws.Range("E7").Value = "  -2.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.84"
$ws.Range("E8").Value = "  +3.88%  "

$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0595"
$ws.Range("E10").Value = "  -0.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0874"
$ws.Range("E11").Value = "  +1.34%  "

$ws.Range("D12").Value = "1.790.23"
$ws.Range("E12").Value = "  -0.08%  "

$ws.Range("D13").Value = "1.582.65"
$ws.Range("E13").Value = "  +0.88%  "

$ws.Range("E14").Value = "  -1.33%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.517"

$ws.Range("D16").Value = "27.681.51"
$ws.Range("E16").Value = "  +1.96%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.18"
$ws.Range("E17").Value = "  +1.63%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.46"
$ws.Range("E18").Value = "  +5.57%  "

$ws.Range("D19").Value = "0.0₃0700"
$ws.Range("E19").Value = "  -0.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.40"
$ws.Range("E20").Value = "  +0.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.988"
$ws.Range("E21").Value = "  -2.12%  "

$ws.Range("E22").Value = "  -0.83%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.28"
$ws.Range("E23").Value = "  +0.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.92"
$ws.Range("E24").Value = "  -0.98%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.18"
$ws.Range("E25").Value = "  -2.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.17"
$ws.Range("E26").Value = "  +0.56%  "

$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.107"
$ws.Range("E27").Value = "  +1.17%  "

$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.55"
$ws.Range("E28").Value = "  -0.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.992"
$ws.Range("E29").Value = "  -1.83%  "

$ws.Range("E30").Value = "  -0.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0469"
$ws.Range("E31").Value = "  -0.66%  "

$ws.Range("E32").Value = "  -0.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.11"
$ws.Range("E33").Value = "  -2.50%  "

$ws.Range("D34").Value = "1.405.80"
$ws.Range("E34").Value = "  -1.72%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("E35").Value = "  -1.90%  "

$ws.Range("E36").Value = "  -4.34%  "

$ws.Range("E37").Value = "  -2.39%  "

$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("E39").Value = "  +2.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.44"
$ws.Range("E40").Value = "  +4.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.804"
$ws.Range("E41").Value = "  -0.47%  "

$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.989"
$ws.Range("E42").Value = "  -2.12%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.61"
$ws.Range("E43").Value = "  -3.88%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.83"
$ws.Range("E44").Value = "  +5.65%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.967"
$ws.Range("E45").Value = "  -3.49%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.45"
$ws.Range("E46").Value = "  -1.62%  "

$ws.Range("D47").Value = "1.700.83"
$ws.Range("E47").Value = "  -0.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.23"
$ws.Range("E48").Value = "  +0.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0523"
$ws.Range("E49").Value = "  +1.03%  "

$ws.Range("D50").Value = "0.0₆0100"
$ws.Range("E50").Value = "  -1.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "39.36"
$ws.Range("E51").Value = "  +15.94%  "

